$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H2").Value = "TestFont"
$ws.Range("H2").Font.OutlineFont = $false
